$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 421.2857
$ws.Range("I4").Value = 416.66666
$ws.Range("K4").Value = 416.66666
$ws.Range("M4").Value = -302.66666
# Row 63
$ws.Range("H63").Value = 29425
$ws.Range("J63").Value = 29425
$ws.Range("L63").Value = 29425
$ws.Range("N63").Value = -30673
# Row 66
$ws.Range("H66").Value = 29425
$ws.Range("J66").Value = 29425
$ws.Range("L66").Value = 88275
$ws.Range("N66").Value = -94515
# Row 75
$ws.Range("H75").Value = 25347
$ws.Range("J75").Value = 25347
$ws.Range("L75").Value = 25347
$ws.Range("N75").Value = -27219
# Row 78
$ws.Range("H78").Value = 25347
$ws.Range("J78").Value = 25347
$ws.Range("L78").Value = 76041
$ws.Range("N78").Value = -85401
# Row 95
$ws.Range("H95").Value = 26425
$ws.Range("J95").Value = 26425
$ws.Range("L95").Value = 26425
$ws.Range("N95").Value = -31917
# Row 112
$ws.Range("H112").Value = 1740
$ws.Range("J112").Value = 1255
$ws.Range("L112").Value = 3765
$ws.Range("N112").Value = -5981
# Row 130
$ws.Range("H130").Value = 30468
$ws.Range("J130").Value = 30468
$ws.Range("L130").Value = 30468
$ws.Range("N130").Value = -40508

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18073.434
$ws.Range("I32").Value = 14995.131
$ws.Range("J32").Value = 28187.857
$ws.Range("K32").Value = 14995.131
$ws.Range("L32").Value = 28187.857
$ws.Range("M32").Value = -14708.131
$ws.Range("N32").Value = -28761.857

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2487.2
$ws.Range("I20").Value = 2330.8
$ws.Range("J20").Value = 2800
$ws.Range("K20").Value = 2330.8
$ws.Range("L20").Value = 2800
$ws.Range("M20").Value = -2083.8
$ws.Range("N20").Value = -3294

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6809.7646
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6809.7646
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6809.7646
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -7399.7646
# Row 34
$ws.Range("H34").Value = 6809.7646
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6809.7646
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6809.7646
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -7213.7646
# Row 52
$ws.Range("H52").Value = 51166.668
$ws.Range("J52").Value = 51166.668
$ws.Range("L52").Value = 51166.668
$ws.Range("N52").Value = -51754.668

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 7053.4443
$ws.Range("J34").Value = 9658.615
$ws.Range("L34").Value = 28975.845
$ws.Range("N34").Value = -29143.845
# Row 39
$ws.Range("H39").Value = 2760
$ws.Range("J39").Value = 2760
$ws.Range("L39").Value = 8280
$ws.Range("N39").Value = -8868
# Row 55
$ws.Range("H55").Value = 2204.4546
$ws.Range("I55").Value = 670
$ws.Range("J55").Value = 3483.1667
$ws.Range("K55").Value = 2010
$ws.Range("L55").Value = 10449.5001
$ws.Range("M55").Value = -1833
$ws.Range("N55").Value = -10803.5001
# Row 59
$ws.Range("H59").Value = 3333.3333
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 3333.3333
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 9999.999899999999
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -11079.9999
# Row 80
$ws.Range("H80").Value = 2800.2222
$ws.Range("I80").Value = 2234
$ws.Range("J80").Value = 3083.3333
$ws.Range("K80").Value = 6702
$ws.Range("L80").Value = 9249.999899999999
$ws.Range("M80").Value = -5766
$ws.Range("N80").Value = -11121.9999
# Row 83
$ws.Range("H83").Value = 2800.2222
$ws.Range("I83").Value = 2234
$ws.Range("J83").Value = 3083.3333
$ws.Range("K83").Value = 20106
$ws.Range("L83").Value = 27749.9997
$ws.Range("M83").Value = -15426
$ws.Range("N83").Value = -37109.9997
# Row 92
$ws.Range("H92").Value = 2150.1
$ws.Range("J92").Value = 2289
$ws.Range("L92").Value = 6867
$ws.Range("N92").Value = -9363
# Row 121
$ws.Range("H121").Value = 15174
$ws.Range("I121").Value = 243.33333
$ws.Range("J121").Value = 35081.555
$ws.Range("K121").Value = 729.99999
$ws.Range("L121").Value = 105244.665
$ws.Range("M121").Value = 580.00001
$ws.Range("N121").Value = -107864.665
# Row 137
$ws.Range("H137").Value = 3208.04
$ws.Range("I137").Value = 2579.375
$ws.Range("K137").Value = 7738.125
$ws.Range("M137").Value = -2638.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 83.875
$ws.Range("I2").Value = 97.25
$ws.Range("K2").Value = 97.25
$ws.Range("M2").Value = 15.75
# Row 57
$ws.Range("H57").Value = 15553.667
$ws.Range("J57").Value = 23661
$ws.Range("L57").Value = 23661
$ws.Range("N57").Value = -25301
# Row 80
$ws.Range("H80").Value = 3545.7144
$ws.Range("I80").Value = 3528.3333
$ws.Range("J80").Value = 3650
$ws.Range("K80").Value = 3528.3333
$ws.Range("L80").Value = 3650
$ws.Range("M80").Value = -2530.3333
$ws.Range("N80").Value = -5646
# Row 83
$ws.Range("H83").Value = 3545.7144
$ws.Range("I83").Value = 3528.3333
$ws.Range("J83").Value = 3650
$ws.Range("K83").Value = 17641.6665
$ws.Range("L83").Value = 18250
$ws.Range("M83").Value = -12649.6665
$ws.Range("N83").Value = -28234
# Row 97
$ws.Range("H97").Value = 1280.9546
$ws.Range("I97").Value = 844.7059
$ws.Range("J97").Value = 2764.2
$ws.Range("K97").Value = 844.7059
$ws.Range("L97").Value = 2764.2
$ws.Range("M97").Value = -348.7059
$ws.Range("N97").Value = -3756.2
# Row 126
$ws.Range("H126").Value = 718643.9
$ws.Range("I126").Value = 2750
$ws.Range("J126").Value = 1005001.4
$ws.Range("K126").Value = 8250
$ws.Range("L126").Value = 3015004.2
$ws.Range("M126").Value = -5780
$ws.Range("N126").Value = -3019944.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 3497.7856
$ws.Range("I100").Value = 1652.7142
$ws.Range("J100").Value = 5342.857
$ws.Range("K100").Value = 1652.7142
$ws.Range("L100").Value = 5342.857
$ws.Range("M100").Value = -1111.7142
$ws.Range("N100").Value = -6424.857
# Row 132
$ws.Range("H132").Value = 2489.6
$ws.Range("I132").Value = 1686.4517
$ws.Range("K132").Value = 5059.355100000001
$ws.Range("M132").Value = -2529.355100000001
# Row 136
$ws.Range("H136").Value = 3454381.8
$ws.Range("I136").Value = 7698261.5
$ws.Range("J136").Value = 6229.375
$ws.Range("K136").Value = 23094784.5
$ws.Range("L136").Value = 18688.125
$ws.Range("M136").Value = -23092234.5
$ws.Range("N136").Value = -23788.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 68366
$ws.Range("J39").Value = 68366
$ws.Range("L39").Value = 68366
$ws.Range("N39").Value = -69192
# Row 42
$ws.Range("H42").Value = 60049
$ws.Range("J42").Value = 60049
$ws.Range("L42").Value = 60049
$ws.Range("N42").Value = -60805
# Row 122
$ws.Range("H122").Value = 2691.2222
$ws.Range("I122").Value = 1788.8572
$ws.Range("J122").Value = 5849.5
$ws.Range("K122").Value = 5366.571599999999
$ws.Range("L122").Value = 17548.5
$ws.Range("M122").Value = -2916.571599999999
$ws.Range("N122").Value = -22448.5
# Row 132
$ws.Range("H132").Value = 1828559
$ws.Range("I132").Value = 2130545.5
$ws.Range("K132").Value = 6391636.5
$ws.Range("M132").Value = -6389106.5
